# LMS-2340 Changes to OD600 and Transcriptomics templates.
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("openbis-metadata")
$wsData = $wb.Worksheets.Item("openbis-data")

# Update the Strain value (B3) from "strain1" to "MGP90"
$wsMeta.Range("B3").Value = "MGP90"

# Add a new "Header Format" metadata row describing the expected
# transcriptomics/proteomics header format, matching the formatting of
# the preceding row.
[void]$wsMeta.Range("A7:C7").Copy()
[void]$wsMeta.Range("A8:C8").PasteSpecial(-4122)
$wsMeta.Range("A8").Value = "Header Format"
$wsMeta.Range("B8").Value = "TIME::VALUE_TYPE"
$wsMeta.Range("C8").Value = "Must be TIME::VALUE_TYPE"

# Widen column B on the metadata sheet so the longer values fit.
$wsMeta.Columns.Item(2).ColumnWidth = 23

# Set the metadata sheet's print page setup (B5 paper, portrait).
$wsMeta.PageSetup.PaperSize = 10
$wsMeta.PageSetup.Orientation = 1

# Make the metadata sheet the active / selected sheet & cell, with the
# data sheet no longer the tab-selected one.
[void]$wsMeta.Select()
[void]$wsMeta.Range("B9").Select()

$wb.Save()
